$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write data rows first (2..35), then the header row (1) last, so the
# shared-string table ends up in the same append order Excel produced
# when the author typed the new KEY/ES/EN header after already having
# the data underneath it.

$ws.Range("A2").Value = "app.name"
$ws.Range("B2").Value = "Gestión de Corpus Documentales"
$ws.Range("C2").Value = "Documental Corpus Management"

$ws.Range("A3").Value = "app.languageSelector.language"
$ws.Range("B3").Value = "Idioma"
$ws.Range("C3").Value = "Language"

$ws.Range("A4").Value = "app.navBar.home"
$ws.Range("B4").Value = "Inicio"
$ws.Range("C4").Value = "Home"

$ws.Range("A5").Value = "app.navBar.addDocuments"
$ws.Range("B5").Value = "Añadir documentos"
$ws.Range("C5").Value = "Add documents"

$ws.Range("A6").Value = "app.navBar.visualizeData"
$ws.Range("B6").Value = "Visualizar datos"
$ws.Range("C6").Value = "Visualize data"

$ws.Range("A7").Value = "menu.intro.title"
$ws.Range("B7").Value = "Bienvenido al servicio web de visualización gráfica de conjuntos de datos"
$ws.Range("C7").Value = "Welcome to your dataset graphic visualization web service."

$ws.Range("A8").Value = "menu.intro.first_paragraph"
$ws.Range("B8").Value = "El mundo de la información evoluciona. Nos vemos obligados cada día a adaptar nuestra tecnología a la oleada de cantidades de datos cada vez más extensas. A la vez que aumenta la necesidad de generar sistemas que aseguren su protección y capacidad de salvaguarda, es necesario el desarrollo de nuevos servicios que nos permitan conocer de una manera auxiliar y sencilla las características de los conjuntos de información que poseemos."
$ws.Range("C8").Value = "The world of global information is evolving. Everyday, we're forced to adapt our technological environment for upcoming waves of extensive data amounts. At the same time it's necessary to generate system to protect information, it's necessary to develop new services to access our documents' most important characteristics, in the simplest of ways."

$ws.Range("A9").Value = "menu.intro.second_paragraph"
$ws.Range("B9").Value = "Por ello, este proyecto, desarrollado inicialmente como un Trabajo de Fin de Grado, pretende poner de manifiesto las utilidades llevadas por tecnologías web que nos permiten analizar nuestros propios corpus documentales, por medio de gráficas interactivas."
$ws.Range("C9").Value = "Because of that, this project - developed initially as a final university degree project - tries to manifest the ultimate web techonologies features which permit us to analyze our own documental corpuses using interactive graphs."

$ws.Range("A10").Value = "menu.indexList.title"
$ws.Range("B10").Value = "Consulta la lista actual de conjuntos en línea"
$ws.Range("C10").Value = "Check the current online index list"

$ws.Range("A11").Value = "menu.indexList.subtitle"
$ws.Range("B11").Value = "Selecciona un corpus para empezar"
$ws.Range("C11").Value = "Select a corpus to begin"

$ws.Range("A12").Value = "menu.indexList.inputPlaceholder"
$ws.Range("B12").Value = "… o crea un nuevo índice"
$ws.Range("C12").Value = "… or create a new one"

$ws.Range("A13").Value = "menu.indexList.create"
$ws.Range("B13").Value = "Crear"
$ws.Range("C13").Value = "Create"

$ws.Range("A14").Value = "menu.indexList.documents"
$ws.Range("B14").Value = "documentos"
$ws.Range("C14").Value = "documents"

$ws.Range("A15").Value = "menu.indexList.creating"
$ws.Range("B15").Value = "Creando el índice…"
$ws.Range("C15").Value = "Creating index…"

$ws.Range("A16").Value = "upload.title"
$ws.Range("B16").Value = "Gestiona los documentos que componen tu corpus desde esta página"
$ws.Range("C16").Value = "Manage the documents from your dataset in this page"

$ws.Range("A17").Value = "upload.elasticConnection"
$ws.Range("B17").Value = "Conexión al servidor de Elasticsearch"
$ws.Range("C17").Value = "Elasticsearch server connection"

$ws.Range("A18").Value = "upload.uploader.addInfo"
$ws.Range("B18").Value = "Añade información"
$ws.Range("C18").Value = "Add information"

$ws.Range("A19").Value = "upload.uploader.uploadDocuments"
$ws.Range("B19").Value = "Sube documentos al corpus"
$ws.Range("C19").Value = "Upload documents to corpus"

$ws.Range("A20").Value = "upload.uploader.drag"
$ws.Range("B20").Value = "Arrastra los documentos aquí"
$ws.Range("C20").Value = "Drag your documents here"

$ws.Range("A21").Value = "upload.uploader.choose"
$ws.Range("B21").Value = "Elige los ficheros"
$ws.Range("C21").Value = "Choose files"

$ws.Range("A22").Value = "upload.uploader.queue"
$ws.Range("B22").Value = "Cola de almacenamiento"
$ws.Range("C22").Value = "Queue"

$ws.Range("A23").Value = "upload.uploader.stillNoFiles"
$ws.Range("B23").Value = "Todavía nmo se han añadido documentos."
$ws.Range("C23").Value = "Still no added files."

$ws.Range("A24").Value = "common.name"
$ws.Range("B24").Value = "Nombre"
$ws.Range("C24").Value = "Name"

$ws.Range("A25").Value = "common.format"
$ws.Range("B25").Value = "Formato"
$ws.Range("C25").Value = "Extension"

$ws.Range("A26").Value = "common.options"
$ws.Range("B26").Value = "Opciones"
$ws.Range("C26").Value = "Options"

$ws.Range("A27").Value = "upload.uploader.remove"
$ws.Range("B27").Value = "Quitar de la cola"
$ws.Range("C27").Value = "Remove from queue"

$ws.Range("A28").Value = "upload.uploader.add"
$ws.Range("B28").Value = "Añadir al corpus"
$ws.Range("C28").Value = "Add to corpus"

$ws.Range("A29").Value = "upload.uploader.addAll"
$ws.Range("B29").Value = "Añadir todo"
$ws.Range("C29").Value = "Add all"

$ws.Range("A30").Value = "upload.uploader.emptyCorpus"
$ws.Range("B30").Value = "El corpus está vacío."
$ws.Range("C30").Value = "Corpus is empty"

$ws.Range("A31").Value = "upload.uploader.cleanCorpus"
$ws.Range("B31").Value = "Borrar datos del corpus"
$ws.Range("C31").Value = "Clean corpus"

$ws.Range("A32").Value = "stats.title"
$ws.Range("B32").Value = "Consulta la información sobre el índice"
$ws.Range("C32").Value = "Check information about index"

$ws.Range("A33").Value = "stats.filterBar.filters"
$ws.Range("B33").Value = "Filtros"
$ws.Range("C33").Value = "Filters"

$ws.Range("A34").Value = "stats.filterBar.noFiltersAvailable"
$ws.Range("B34").Value = "Por el momento no hay filtros aplicables."
$ws.Range("C34").Value = "No active filters for the moment."

$ws.Range("A35").Value = "stats.searchBar.search"
$ws.Range("B35").Value = "Busca términos dentro del corpus"
$ws.Range("C35").Value = "Search terms inside the corpus"

$ws.Range("A1").Value = "KEY"
$ws.Range("B1").Value = "ES"
$ws.Range("C1").Value = "EN"

# Match the saved selection shown in the diff ($ws.Range("B15") was
# the active cell when the author saved).
$ws.Range("B15").Select() | Out-Null
